$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This is the "Generate Report for Handback" edit: the handback run stamped
# each localized-file row with its target markdown link, the handback xliff
# file name, and the handback timestamp, and flipped the overview "Ready for
# handoff" status to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------

$colWidth30 = 29.16666666666667   # -> stored col width ~29.98 (status cols)
$colWidth40 = 39.16666666666667   # -> stored col width 40    (target/handback cols)

$statusText = "Handed back: in sync with en-US"

# ----------------------- Overview sheet ------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText
$ov.Columns.Item(5).ColumnWidth = $colWidth30
$ov.Columns.Item(6).ColumnWidth = $colWidth30

# ----------------------- zh-cn sheet ----------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText
$zh.Columns.Item(3).ColumnWidth = $colWidth30
$zh.Columns.Item(9).ColumnWidth = $colWidth40
$zh.Columns.Item(10).ColumnWidth = $colWidth40

$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f897ecb65853a8dcd97e3447b42525b45618efdb/e2e/3c3d9d13-1f96-425b-b80d-1a40828f183c.md", "", "", "3c3d9d13-1f96-425b-b80d-1a40828f183c.md")
$zh.Range("J2").Value = "3c3d9d13-1f96-425b-b80d-1a40828f183c.3bc0c8fa02ec56afce51dc431dfc244e4cc4d249.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-26 10:59:28"

$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f897ecb65853a8dcd97e3447b42525b45618efdb/e2e/9d9e1ce4-379c-4b62-b206-f5d5b2ff08cc.md", "", "", "9d9e1ce4-379c-4b62-b206-f5d5b2ff08cc.md")
$zh.Range("J3").Value = "9d9e1ce4-379c-4b62-b206-f5d5b2ff08cc.114a0528c0fdda69b6bdb0317dfd4520918ae185.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-26 10:59:28"

# ----------------------- de-de sheet ----------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText
$de.Columns.Item(3).ColumnWidth = $colWidth30
$de.Columns.Item(9).ColumnWidth = $colWidth40
$de.Columns.Item(10).ColumnWidth = $colWidth40

$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f897ecb65853a8dcd97e3447b42525b45618efdb/e2e/3c3d9d13-1f96-425b-b80d-1a40828f183c.md", "", "", "3c3d9d13-1f96-425b-b80d-1a40828f183c.md")
$de.Range("J2").Value = "3c3d9d13-1f96-425b-b80d-1a40828f183c.3bc0c8fa02ec56afce51dc431dfc244e4cc4d249.de-de.xlf"
$de.Range("K2").Value = "2016-08-26 10:59:34"

$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f897ecb65853a8dcd97e3447b42525b45618efdb/e2e/9d9e1ce4-379c-4b62-b206-f5d5b2ff08cc.md", "", "", "9d9e1ce4-379c-4b62-b206-f5d5b2ff08cc.md")
$de.Range("J3").Value = "9d9e1ce4-379c-4b62-b206-f5d5b2ff08cc.114a0528c0fdda69b6bdb0317dfd4520918ae185.de-de.xlf"
$de.Range("K3").Value = "2016-08-26 10:59:34"
